# Loan RBI, Variable Instalments
# Insert a new column before column N on the "Repayment Schedule" sheet,
# splitting the previous "In Advance"/data column away from the "Due" column,
# and make this sheet the active tab/selection (previously "Input" was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Range("K19").Select()
